$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry data for the affected rows (A..AY minus the
# handful of always-blank placeholder columns, which we leave untouched so
# we don't disturb cells that are empty on both sides of the swap).
$dataCols = @(1,2,4,5,6,7,8,16,17,18,19,20,21,22,23,25,26,27,28,29,30,31,33,36,37,41,49,50)

function Get-RowSnapshot($ws, $row, $cols) {
    $snap = @{}
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($row, $c)
        $snap[$c] = $cell.Value2
    }
    return $snap
}

function Set-RowFromSnapshot($ws, $row, $snap, $cols) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($row, $c)
        $v = $snap[$c]
        if ($v -eq $null) {
            $cell.ClearContents()
        } elseif ($v -is [string]) {
            if ($v -eq "") {
                $cell.ClearContents()
            } else {
                # Force text so Excel doesn't auto-convert date/time-looking
                # strings (e.g. "2026-01-20", "13:20") into date serials.
                $cell.NumberFormat = "@"
                $cell.Value = $v
                $cell.Style = "Normal"
            }
        } else {
            $cell.Value = $v
        }
    }
}

# Rows participating in the reshuffle, snapshotted BEFORE any writes since
# some rows are both a source and a destination.
$rowsToSnapshot = @(46, 47, 54, 55, 56, 57, 58, 60, 61)
$snaps = @{}
foreach ($r in $rowsToSnapshot) {
    $snaps[$r] = Get-RowSnapshot $ws $r $dataCols
}

# destination row -> source row (content that ends up there)
$mapping = @{
    46 = 47
    47 = 46
    54 = 55
    55 = 54
    56 = 57
    57 = 58
    58 = 56
    60 = 61
    61 = 60
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    Set-RowFromSnapshot $ws $dst $snaps[$src] $dataCols
}

# Row 55 picks up row 54's content, which (unlike the current row 55) has no
# J/K/N/AF entries, so those leftover cells must be cleared out explicitly.
foreach ($col in @("J", "K", "N", "AF")) {
    $ws.Range($col + "55").ClearContents()
}
